$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record arrived for Achicoria (Vega Central Mapocho de
# Santiago). Insert a fresh row right above the existing row 45 so the
# later rows (old 45, 46, 47) shift down to 46, 47, 48 unchanged, then
# populate the newly inserted row 45 with the new record.
$ws.Rows.Item(45).Insert()

$ws.Range("A45").Value = 9
$ws.Range("B45").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C45").Value = "Metropolitana"
$ws.Range("D45").Value = 45013
$ws.Range("E45").Value = 13
$ws.Range("F45").Value = 100112010
$ws.Range("G45").Value = "Achicoria"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 70
$ws.Range("K45").Value = 7000
$ws.Range("L45").Value = 7000
$ws.Range("M45").Value = 7000
$ws.Range("N45").Value = "$/caja 16 unidades"
$ws.Range("O45").Value = "Provincia de Quillota"
$ws.Range("P45").Value = 438
$ws.Range("Q45").Value = 16
$ws.Range("R45").Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of
# column D (style index 2 in the original workbook).
$ws.Range("D45").NumberFormat = $ws.Range("D44").NumberFormat
